# Remove the "O/O Prevent submarines from leaving port" bullet paragraph
# from the Main Effort / additional target effects list, per the commit:
# "Updated additional target effects."
$d = $word.ActiveDocument

$target = "O/O Prevent submarines from leaving port"
$removed = $false

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    if ($text -ne $null -and $text.TrimEnd([char]13, [char]7) -eq $target) {
        $p.Range.Delete()
        $removed = $true
        break
    }
}

if (-not $removed) {
    throw "Target paragraph not found: $target"
}

Write-Output "removed=$removed"
